# Updated research Data sheet
# Applies the data corrections made to the SearchData and UpdatedData sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SearchData sheet (sheet1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SearchData")

$searchDataChanges = @{
    2   = @{ N = 45;  Q = 75;  R = 70;  V = 13; W = 160; X = 160; Y = 160 }
    3   = @{ N = 45;  Q = 75;  R = 70;  V = 13; W = 160; X = 160; Y = 160 }
    4   = @{ N = 45;  Q = 75;  R = 70;  V = 13; W = 160; X = 160; Y = 160 }
    6   = @{ Q = 35 }
    7   = @{ N = 13;  Q = 7;   R = 7;   X = 18 }
    8   = @{ N = 25;  R = 8;   X = 48 }
    14  = @{ N = 65 }
    16  = @{ N = 35 }
    18  = @{ N = 35 }
    23  = @{ W = 30 }
    24  = @{ Q = 2 }
    25  = @{ Q = 2 }
    26  = @{ Q = 2 }
    30  = @{ N = 65;  Q = 110; R = 110 }
    31  = @{ N = 65;  Q = 110; R = 110 }
    32  = @{ N = 65;  Q = 110; R = 110 }
    33  = @{ N = 65;  Q = 110; R = 110 }
    40  = @{ N = 5;   Y = 17 }
    41  = @{ N = 45;  V = 11;  Y = 160 }
    49  = @{ Q = 96;  R = 96 }
    50  = @{ N = 45 }
    78  = @{ N = 45 }
    79  = @{ Q = 40 }
    82  = @{ Q = 25 }
    96  = @{ Q = 2 }
    97  = @{ Q = 2 }
    98  = @{ Q = 2 }
    106 = @{ X = 3 }
    116 = @{ S = 10 }
    123 = @{ S = 10;  T = 1;   U = 1;   V = 11 }
    126 = @{ Q = 1 }
}

foreach ($row in $searchDataChanges.Keys) {
    $cols = $searchDataChanges[$row]
    foreach ($col in $cols.Keys) {
        $ws1.Range("$col$row").Value = $cols[$col]
    }
}

# Cells that also pick up formatting used elsewhere in the sheet
# X4 should match the shaded/bordered look already used by X2/X3
$ws1.Range("X2").Copy()
$ws1.Range("X4").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# X106 and V123 pick up a plain yellow highlight fill
$ws1.Range("X106").Interior.Color = 65535
$ws1.Range("V123").Interior.Color = 65535

# View state: scrolled/selected cell moved
$ws1.Activate()
$ws1.Range("C2").Select()

# ---------------------------------------------------------------------------
# UpdatedData sheet (sheet2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("UpdatedData")

$ws2.Range("AF9").Value = 160
$ws2.Range("AG9").Value = 160
$ws2.Range("AH9").Value = 160

$ws2.Activate()
$ws2.Range("A9").Select()
